$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "button_closeAlert_class"
$ws.Range("B1").Value = "button_closeAlert_class_1"
$ws.Range("C1").Value = "div_backdropMenu_class"
$ws.Range("D1").Value = "div_backdropMenu_class_1"
$ws.Range("F1").Value = "p_sessionInfo_class"
